# Generacion de metodos POST modificados
# Duplicate the "post" sheet into a new "post2" sheet, then drop the
# "entrada / valor1" input row so post2 only keeps the BODY rows and the
# salida/gastoTotal row.

$wb = $excel.ActiveWorkbook

$postSheet = $wb.Worksheets.Item("post")

# Place the copy right after "post" (i.e. at the end of the workbook).
$postSheet.Copy($null, $postSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "post2"

# Remove row 2 ("entrada" / "valor1" / 5.0) - post2 starts straight at the
# body rows (row 2 becomes the old row 3, etc.)
$newSheet.Rows.Item(2).Delete()
